$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update StructureDefinition metadata properties ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: refresh publication timestamp
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was previously empty; now populated
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 previously held "Contact" / "No display for ContactDetail".
# It becomes "Jurisdiction" / "United States of America".
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row - remove it,
# shifting all subsequent rows up by one (A1:B21 -> A1:B20).
$meta.Rows.Item(11).Delete()

# --- Sheet "Elements": fix the root Extension row's Short/Definition text ---
$elements = $wb.Worksheets.Item("Elements")

# K2 (Short): "Extension" -> "Claim Snapshot Provider Name"
$elements.Range("K2").Value = "Claim Snapshot Provider Name"

# L2 (Definition): "An Extension" -> "Original provider name as reported on the claim"
$elements.Range("L2").Value = "Original provider name as reported on the claim"
